$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new data rows before the current row 260 (Artic Snow / 44293 / ...),
# pushing the old rows 260-268 down to 263-271.
$ws.Rows("260:262").Insert()

# --- New row 260: Ruby Diamond, Primera ---
$ws.Range("A260").Value = 7
$ws.Range("B260").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C260").Value = "Ñuble"
$ws.Range("D260").Value = 44578
$ws.Range("E260").Value = 16
$ws.Range("F260").Value = "Fruta"
$ws.Range("G260").Value = 100103
$ws.Range("H260").Value = "Frutos de hueso (carozo)"
$ws.Range("I260").Value = 100103006
$ws.Range("J260").Value = "Nectarín"
$ws.Range("K260").Value = "Ruby Diamond"
$ws.Range("L260").Value = "Primera"
$ws.Range("M260").Value = 300
$ws.Range("N260").Value = 13000
$ws.Range("O260").Value = 14000
$ws.Range("P260").Value = 13500
$ws.Range("Q260").Value = "$/caja 16 kilos empedrada"
$ws.Range("R260").Value = "Región de O'Higgins"
$ws.Range("S260").Value = 844
$ws.Range("T260").Value = 16

# --- New row 261: Ruby Diamond, Segunda ---
$ws.Range("A261").Value = 7
$ws.Range("B261").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C261").Value = "Ñuble"
$ws.Range("D261").Value = 44578
$ws.Range("E261").Value = 16
$ws.Range("F261").Value = "Fruta"
$ws.Range("G261").Value = 100103
$ws.Range("H261").Value = "Frutos de hueso (carozo)"
$ws.Range("I261").Value = 100103006
$ws.Range("J261").Value = "Nectarín"
$ws.Range("K261").Value = "Ruby Diamond"
$ws.Range("L261").Value = "Segunda"
$ws.Range("M261").Value = 120
$ws.Range("N261").Value = 11000
$ws.Range("O261").Value = 12000
$ws.Range("P261").Value = 11500
$ws.Range("Q261").Value = "$/caja 16 kilos empedrada"
$ws.Range("R261").Value = "Región de O'Higgins"
$ws.Range("S261").Value = 719
$ws.Range("T261").Value = 16

# --- New row 262: Super Queen, Primera ---
$ws.Range("A262").Value = 7
$ws.Range("B262").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C262").Value = "Ñuble"
$ws.Range("D262").Value = 44578
$ws.Range("E262").Value = 16
$ws.Range("F262").Value = "Fruta"
$ws.Range("G262").Value = 100103
$ws.Range("H262").Value = "Frutos de hueso (carozo)"
$ws.Range("I262").Value = 100103006
$ws.Range("J262").Value = "Nectarín"
$ws.Range("K262").Value = "Super Queen"
$ws.Range("L262").Value = "Primera"
$ws.Range("M262").Value = 160
$ws.Range("N262").Value = 11000
$ws.Range("O262").Value = 12000
$ws.Range("P262").Value = 11500
$ws.Range("Q262").Value = "$/caja 16 kilos empedrada"
$ws.Range("R262").Value = "Región de O'Higgins"
$ws.Range("S262").Value = 719
$ws.Range("T262").Value = 16
